# Natmi following Dr Hou advice
# Update ligand/receptor-expressing-cell counts and the resulting
# derived statistic columns for each Wnt2-Fzd9 edge row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3884013333333334
$ws.Range("H2").Value = 1.165204
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3962276666666666
$ws.Range("N2").Value = 1.188683
$ws.Range("O2").Value = 0.1703050239984229
$ws.Range("P2").Value = 0.1960950025990754
$ws.Range("Q2").Value = 0.1538953540368889
$ws.Range("R2").Value = 1.385058186332
$ws.Range("S2").Value = 0.1703050239984229
$ws.Range("T2").Value = 0.1960950025990754

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3884013333333334
$ws.Range("H3").Value = 1.165204
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6295006666666667
$ws.Range("N3").Value = 1.888502
$ws.Range("O3").Value = 0.270569511325618
$ws.Range("P3").Value = 0.3115429467724861
$ws.Range("Q3").Value = 0.2444988982675556
$ws.Range("R3").Value = 2.200490084408
$ws.Range("S3").Value = 0.270569511325618
$ws.Range("T3").Value = 0.3115429467724861

# Row 4 (Target cluster: M1)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3884013333333334
$ws.Range("H4").Value = 1.165204
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1632393333333333
$ws.Range("N4").Value = 0.489718
$ws.Range("O4").Value = 0.07016289098309612
$ws.Range("P4").Value = 0.08078794134585419
$ws.Range("Q4").Value = 0.06340237471911112
$ws.Range("R4").Value = 0.5706213724720001
$ws.Range("S4").Value = 0.07016289098309612
$ws.Range("T4").Value = 0.08078794134585419

# Row 5 (Target cluster: M2)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3884013333333334
$ws.Range("H5").Value = 1.165204
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2196503333333334
$ws.Range("N5").Value = 0.6589510000000001
$ws.Range("O5").Value = 0.0944092460889781
$ws.Range("P5").Value = 0.1087060200723518
$ws.Range("Q5").Value = 0.08531248233377781
$ws.Range("R5").Value = 0.7678123410040002
$ws.Range("S5").Value = 0.0944092460889781
$ws.Range("T5").Value = 0.1087060200723518

# Row 6 (Target cluster: sCs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3884013333333334
$ws.Range("H6").Value = 1.165204
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 0.9179585
$ws.Range("N6").Value = 1.835917
$ws.Range("O6").Value = 0.394553327603885
$ws.Range("P6").Value = 0.3028680892102324
$ws.Range("Q6").Value = 0.3565363053446667
$ws.Range("R6").Value = 2.139217832068
$ws.Range("S6").Value = 0.394553327603885
$ws.Range("T6").Value = 0.3028680892102324
